$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("scenario_info")

# Remove the "networkCrs" parameter row (row 6) - CRS is now taken directly from the file
$ws.Rows.Item(6).Select()
$ws.Rows.Item(6).Delete()

# Re-apply the autofilter so its range shrinks along with the data (A1:E25 -> A1:E24)
$ws.AutoFilterMode = $false
$ws.Range("A1:E24").AutoFilter()

# Make sure the hidden _FilterDatabase defined name matches the new filter range
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "=scenario_info!`$A`$1:`$E`$24"

$wb.Save()
